# Update "想去人数" (column F) values across the sheets, matching the
# gh-pages data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 1673
$ws.Cells.Item(6, 6).Value = 621
$ws.Cells.Item(7, 6).Value = 1111
$ws.Cells.Item(8, 6).Value = 1554
$ws.Cells.Item(12, 6).Value = 1464
$ws.Cells.Item(13, 6).Value = 3094
$ws.Cells.Item(14, 6).Value = 630
$ws.Cells.Item(15, 6).Value = 1776
$ws.Cells.Item(16, 6).Value = 1808
$ws.Cells.Item(17, 6).Value = 859
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(20, 6).Value = 1475
$ws.Cells.Item(23, 6).Value = 11
$ws.Cells.Item(24, 6).Value = 1218
$ws.Cells.Item(25, 6).Value = 407
$ws.Cells.Item(26, 6).Value = 459
$ws.Cells.Item(27, 6).Value = 121
$ws.Cells.Item(28, 6).Value = 4804
$ws.Cells.Item(29, 6).Value = 43
$ws.Cells.Item(30, 6).Value = 753
$ws.Cells.Item(32, 6).Value = 1659
$ws.Cells.Item(33, 6).Value = 73
$ws.Cells.Item(34, 6).Value = 147

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 76

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(10, 6).Value = 76
$ws.Cells.Item(12, 6).Value = 1673
$ws.Cells.Item(14, 6).Value = 621
$ws.Cells.Item(15, 6).Value = 1111
$ws.Cells.Item(16, 6).Value = 1554
$ws.Cells.Item(21, 6).Value = 1464
$ws.Cells.Item(22, 6).Value = 3094
$ws.Cells.Item(23, 6).Value = 630
$ws.Cells.Item(24, 6).Value = 1776
$ws.Cells.Item(25, 6).Value = 1808
$ws.Cells.Item(26, 6).Value = 859
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(29, 6).Value = 1475
$ws.Cells.Item(33, 6).Value = 11
$ws.Cells.Item(35, 6).Value = 1218
$ws.Cells.Item(36, 6).Value = 407
$ws.Cells.Item(37, 6).Value = 459
$ws.Cells.Item(38, 6).Value = 122
$ws.Cells.Item(39, 6).Value = 4804
$ws.Cells.Item(40, 6).Value = 43
$ws.Cells.Item(41, 6).Value = 753
$ws.Cells.Item(43, 6).Value = 1659
$ws.Cells.Item(46, 6).Value = 73
$ws.Cells.Item(47, 6).Value = 147
